{"js": "// \"Reverted to version 3\": the title text \"Version 6\" becomes \"Version 3\".\n// The original two runs (\"Version \" + \"6\") are combined into a single run\n// (\"Version 3\") by replacing the text of the whole paragraph in place, so\n// the existing paragraph/run formatting (Arial Black, 36pt) and the\n// _GoBack bookmark are preserved, matching Word's own run-merging\n// behaviour when the replaced text shares one formatting.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  paragraph.load(\"text\");\n}\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (/Version\\s*6/.test(paragraph.text)) {\n    const newText = paragraph.text.replace(/Version\\s*6/, \"Version 3\");\n    const range = paragraph.getRange();\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# \"Reverted to version 3\": the title text \"Version 6\" becomes \"Version 3\".\n# Using Find/Replace across the whole document merges the old two runs\n# (\"Version \" + \"6\") into a single run (\"Version 3\"), preserving the\n# existing paragraph/run formatting (Arial Black, 36pt) and the _GoBack\n# bookmark, matching Word's own run-merging behaviour.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Version 6\"\n$find.Replacement.Text = \"Version 3\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute([ref]$find.Text, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
